$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "2022-Q1" sheet right before "总计", cloning the layout of
#    "2021-Q4" (same column headers/format) and filling in the new quarter's
#    fund-holding data.
# ---------------------------------------------------------------------------
$prevQ = $wb.Worksheets.Item("2021-Q4")
$newQ = $wb.Worksheets.Add($null, $prevQ)
$newQ.Name = "2022-Q1"

# Clone formatting (styles, borders, fonts) from the previous quarter sheet.
$prevQ.Range("A1:H3").Copy($newQ.Range("A1"))
$newQ.Range("A1").Clear()

# Column D's header differs slightly between quarters ("基金规模" here).
$newQ.Range("D1").Value = "基金规模"

# Make sure these columns are written as text (not numeric) so values such
# as "002379" keep their leading zero and "1.84" stays an exact string.
$newQ.Range("B2:B3").NumberFormat = "@"
$newQ.Range("D2:G3").NumberFormat = "@"

$newQ.Range("A2").Value = 0
$newQ.Range("B2").Value = "002379"
$newQ.Range("C2").Value = "工银瑞信香港中小盘股票（QDII）人民币"
$newQ.Range("D2").Value = "1.84"
$newQ.Range("E2").Value = "86.48"
$newQ.Range("F2").Value = "7.43"
$newQ.Range("G2").Value = "0.1367"
$newQ.Range("H2").Value = 1

$newQ.Range("A3").Value = 1
$newQ.Range("B3").Value = "002380"
$newQ.Range("C3").Value = "工银瑞信香港中小盘股票（QDII）美元"
$newQ.Range("D3").Value = "1.84"
$newQ.Range("E3").Value = "86.48"
$newQ.Range("F3").Value = "7.43"
$newQ.Range("G3").Value = "0.1367"
$newQ.Range("H3").Value = 1

# Drop the plain/default formatting back to unstyled cells (matches the
# other quarter sheets, where only the header row and column A are styled).
$newQ.Range("B2:B3").ClearFormats()
$newQ.Range("D2:G3").ClearFormats()

# ---------------------------------------------------------------------------
# 2. Add the corresponding summary row to "总计": insert a new row 2 for
#    2022-Q1, pushing the older quarters down, and renumber the running
#    index in column A.
# ---------------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")
$tot.Rows("2:2").Insert()

$tot.Range("A3:D3").Copy()
$tot.Range("A2:D2").PasteSpecial(-4122)

$tot.Range("A2").Value = 0
$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 2
$tot.Range("D2").Value = 0.27

$tot.Range("A3").Value = 1
$tot.Range("A4").Value = 2
$tot.Range("A5").Value = 3
$tot.Range("A6").Value = 4
$tot.Range("A7").Value = 5

# Restore the original active sheet/selection (sheet creation shifts focus
# to the new sheet by default).
$wb.Worksheets.Item("2020-Q4").Activate()
